$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of an existing header cell (E1) onto the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean values for rows 2..19, columns F (KNN_Outliers_MAD),
# G (SVM_Outliers_MAD), H (RF_Outliers_MAD).
# Only row 2 / column H (RF) is flagged TRUE; everything else is FALSE.
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}

$ws.Cells.Item(2, 8).Value = $true
